# Edit workbook per commit 'New PO forecast model'
$wb = $excel.ActiveWorkbook

$dateFmt = "YYYY-MM-DD HH:MM:SS"

# --- Sheet: Weekly Quantity (add rows 80-82) ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$weeklyNewRows = @(
    @(80, 45669.99999999999, 33),
    @(81, 45676.99999999999, 132),
    @(82, 45683.99999999999, 5)
)
foreach ($row in $weeklyNewRows) {
    $r = $row[0]
    $dateVal = $row[1]
    $qty = $row[2]
    $wsWeekly.Cells.Item($r, 1).NumberFormat = $dateFmt
    $wsWeekly.Cells.Item($r, 1).Value = $dateVal
    $wsWeekly.Cells.Item($r, 2).Value = $qty
}

# --- Sheet: Monthly Trend (add row 23) ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$monthlyNewRows = ,@(
    @(23, 45688.99999999999, 170)
)
foreach ($row in $monthlyNewRows) {
    $r = $row[0]
    $dateVal = $row[1]
    $qty = $row[2]
    $wsMonthly.Cells.Item($r, 1).NumberFormat = $dateFmt
    $wsMonthly.Cells.Item($r, 1).Value = $dateVal
    $wsMonthly.Cells.Item($r, 2).Value = $qty
}

# --- Sheet: PO Forecast (new forecast model) ---
$wsForecast = $wb.Worksheets.Item("PO Forecast")

# Rows 2-79: only the forecast quantity (column B) changes
$forecastUpdates = @(
    @(2, 1103),
    @(3, 318),
    @(4, 0),
    @(5, 0),
    @(6, 344),
    @(7, 558),
    @(8, 360),
    @(9, 53),
    @(10, 73),
    @(11, 494),
    @(12, 950),
    @(13, 193),
    @(14, 0),
    @(15, 137),
    @(16, 160),
    @(17, 52),
    @(18, 170),
    @(19, 408),
    @(20, 571),
    @(21, 572),
    @(22, 426),
    @(23, 239),
    @(24, 146),
    @(25, 201),
    @(26, 337),
    @(27, 439),
    @(28, 461),
    @(29, 455),
    @(30, 462),
    @(31, 423),
    @(32, 257),
    @(33, 21),
    @(34, 0),
    @(35, 234),
    @(36, 780),
    @(37, 1219),
    @(38, 1231),
    @(39, 852),
    @(40, 417),
    @(41, 203),
    @(42, 178),
    @(43, 139),
    @(44, 59),
    @(45, 1329),
    @(46, 559),
    @(47, 0),
    @(48, 0),
    @(49, 374),
    @(50, 656),
    @(51, 521),
    @(52, 198),
    @(53, 139),
    @(54, 506),
    @(55, 992),
    @(56, 1152),
    @(57, 855),
    @(58, 372),
    @(59, 73),
    @(60, 79),
    @(61, 220),
    @(62, 272),
    @(63, 180),
    @(64, 63),
    @(65, 69),
    @(66, 233),
    @(67, 656),
    @(68, 563),
    @(69, 373),
    @(70, 42),
    @(71, 260),
    @(72, 782),
    @(73, 1271),
    @(74, 1368),
    @(75, 1039),
    @(76, 587),
    @(77, 327),
    @(78, 282),
    @(79, 172)
)
foreach ($row in $forecastUpdates) {
    $r = $row[0]
    $qty = $row[1]
    $wsForecast.Cells.Item($r, 2).Value = $qty
}

# Rows 80-87: both date (column A) and forecast quantity (column B) change
$forecastShifted = @(
    @(80, 45669.99999999999, 805),
    @(81, 45676.99999999999, 142),
    @(82, 45683.99999999999, 32),
    @(83, 45690.99999999999, 400),
    @(84, 45697.99999999999, 738),
    @(85, 45704.99999999999, 676),
    @(86, 45711.99999999999, 352),
    @(87, 45718.99999999999, 220)
)
foreach ($row in $forecastShifted) {
    $r = $row[0]
    $dateVal = $row[1]
    $qty = $row[2]
    $wsForecast.Cells.Item($r, 1).Value = $dateVal
    $wsForecast.Cells.Item($r, 2).Value = $qty
}

# Rows 88-90: brand new rows appended to the forecast
$forecastNewRows = @(
    @(88, 45725.99999999999, 522),
    @(89, 45732.99999999999, 1022),
    @(90, 45739.99999999999, 1264)
)
foreach ($row in $forecastNewRows) {
    $r = $row[0]
    $dateVal = $row[1]
    $qty = $row[2]
    $wsForecast.Cells.Item($r, 1).NumberFormat = $dateFmt
    $wsForecast.Cells.Item($r, 1).Value = $dateVal
    $wsForecast.Cells.Item($r, 2).Value = $qty
}

